$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Ркацители" row (row 4)
$ws.Rows(4).Delete()

# Delete the "Киндзмараули" row (originally row 8, now row 7 after the shift)
$ws.Rows(7).Delete()

# Clear the contents of the now-trailing row (originally row 10, now row 8)
$ws.Range("A8:F8").ClearContents()

# Select A9 per the saved view state
$ws.Range("A9").Select()
